# Refresh the cryptocurrency price/volume table (GitHub Actions scheduled update).
# Price (column D) and Volume(1h) (column E) cells hold plain text values
# (e.g. "96.669.12", "  -0.81%  "). Some new prices look like ordinary
# decimal numbers (e.g. "241.01"); a leading apostrophe forces Excel to
# keep them as text instead of auto-converting to a numeric value, which
# mirrors how the cells were already stored (t="inlineStr") in the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.669.12'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('D3').Value = '3.676.73'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''241.01'
$ws.Range('E5').Value = '  -1.12%  '
$ws.Range('D6').Value = '''1.86'
$ws.Range('E6').Value = '  +9.62%  '
$ws.Range('D7').Value = '''661.66'
$ws.Range('E7').Value = '  +0.47%  '
$ws.Range('D8').Value = '''0.422'
$ws.Range('E8').Value = '  +0.66%  '
$ws.Range('E9').Value = '  +1.98%  '
$ws.Range('E10').Value = '  +0.08%  '
$ws.Range('D11').Value = '3.674.31'
$ws.Range('E11').Value = '  +1.73%  '
$ws.Range('D12').Value = '''45.98'
$ws.Range('E12').Value = '  +4.33%  '
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').Value = '''6.77'
$ws.Range('E14').Value = '  +4.60%  '
$ws.Range('D15').Value = '4.363.69'
$ws.Range('E15').Value = '  +1.90%  '
$ws.Range('E16').Value = '  +3.12%  '
$ws.Range('D17').Value = '96.408.76'
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('D18').Value = '''8.94'
$ws.Range('E18').Value = '  +10.23%  '
$ws.Range('D19').Value = '3.683.57'
$ws.Range('E19').Value = '  +1.98%  '
$ws.Range('D20').Value = '''12.88'
$ws.Range('E20').Value = '  +0.56%  '
$ws.Range('D21').Value = '''18.72'
$ws.Range('E21').Value = '  +3.13%  '
$ws.Range('D22').Value = '''0.527'
$ws.Range('E22').Value = '  -1.45%  '
$ws.Range('D23').Value = '''527.76'
$ws.Range('E23').Value = '  +2.80%  '
$ws.Range('D24').Value = '''3.45'
$ws.Range('E24').Value = '  +0.55%  '
$ws.Range('B25').Value = 'NEARProtocol'
$ws.Range('C25').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D25').Value = '''7.08'
$ws.Range('E25').Value = '  +2.86%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = '''0.0000204'
$ws.Range('E26').Value = '  -1.03%  '
$ws.Range('D27').Value = '''102.38'
$ws.Range('E27').Value = '  +3.59%  '
$ws.Range('D28').Value = '''13.15'
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('D29').Value = '3.876.95'
$ws.Range('E29').Value = '  +1.92%  '
$ws.Range('E30').Value = '  +9.86%  '
$ws.Range('D31').Value = '''12.69'
$ws.Range('E31').Value = '  +7.78%  '
$ws.Range('D32').Value = '''3.05'
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('D33').Value = '''1.00'
$ws.Range('E33').Value = '  +0.19%  '
$ws.Range('D34').Value = '''1.89'
$ws.Range('E34').Value = '  +16.10%  '
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('D36').Value = '''32.65'
$ws.Range('E36').Value = '  +2.58%  '
$ws.Range('D37').Value = '''0.999'
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('D38').Value = '''660.40'
$ws.Range('E38').Value = '  +6.72%  '
$ws.Range('E39').Value = '  +3.30%  '
$ws.Range('D40').Value = '''8.84'
$ws.Range('E40').Value = '  +0.36%  '
$ws.Range('D41').Value = '''43.41'
$ws.Range('E41').Value = '  +30.24%  '
$ws.Range('E42').Value = '  +5.15%  '
$ws.Range('E43').Value = '  +2.82%  '
$ws.Range('D44').Value = '''0.967'
$ws.Range('E44').Value = '  +3.59%  '
$ws.Range('D45').Value = '''6.48'
$ws.Range('E45').Value = '  +8.53%  '
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').Value = '''0.0469'
$ws.Range('E47').Value = '  +6.70%  '
$ws.Range('D48').Value = '''0.453'
$ws.Range('E48').Value = '  +17.86%  '
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').Value = '''23.65'
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').Value = '''8.67'
$ws.Range('E51').Value = '  +1.83%  '
